{"js": "// Change the UML relation line \"ToyBox *-- Toy\" to\n// \"ToyBox \"1\" *-- \"many\" Toy\", wrapping \"ToyBox\" in spell-check proofing\n// marks and moving the \"_GoBack\" bookmark from the \"get_toys()\" paragraph\n// to the end of this paragraph.\n\n// Step 1: remove the existing \"_GoBack\" bookmark from wherever it lives\n// today (the \"get_toys()\" paragraph in before.docx).\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// Step 2: locate the paragraph that still reads \"ToyBox *-- Toy\".\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === \"ToyBox *-- Toy\") {\n    target = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error('Could not find paragraph with text \"ToyBox *-- Toy\"');\n}\n\n// Step 3: rebuild that paragraph's runs: \"ToyBox\" stays wrapped in\n// spellStart/spellEnd proofing marks, followed by new runs for the\n// relationship text, and finish with the relocated bookmark.\nconst ooxml =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>' +\n  '<w:p w:rsidR=\"00C055B4\" w:rsidRDefault=\"00C055B4\" w:rsidP=\"00C055B4\">' +\n  '<w:proofErr w:type=\"spellStart\"/>' +\n  '<w:r><w:t>Toy</w:t></w:r>' +\n  '<w:r w:rsidR=\"00172387\"><w:t>B</w:t></w:r>' +\n  '<w:r><w:t>ox</w:t></w:r>' +\n  '<w:proofErr w:type=\"spellEnd\"/>' +\n  '<w:r><w:t xml:space=\"preserve\"> </w:t></w:r>' +\n  '<w:r><w:t>&quot;1&quot; *-- &quot;many&quot;</w:t></w:r>' +\n  '<w:r><w:t xml:space=\"preserve\"> Toy</w:t></w:r>' +\n  '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/>' +\n  '<w:bookmarkEnd w:id=\"0\"/>' +\n  '</w:p>' +\n  '</w:body></w:document>' +\n  '</pkg:xmlData></pkg:part></pkg:package>';\n\ntarget.insertOoxml(ooxml, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Change the UML relation line \"ToyBox *-- Toy\" to\n# \"ToyBox \"1\" *-- \"many\" Toy\", wrapping \"ToyBox\" in spell-check proofing\n# marks and moving the \"_GoBack\" bookmark from the \"get_toys()\" paragraph\n# to the end of this paragraph.\n\n$d = $word.ActiveDocument\n\n# Step 1: remove the existing \"_GoBack\" bookmark from wherever it lives\n# today (the \"get_toys()\" paragraph in before.docx). \"_GoBack\" is a\n# hidden bookmark, but it can still be addressed directly by name.\n$bm = $d.Bookmarks.Item(\"_GoBack\")\n$bm.Delete()\n\n# Step 2: locate the paragraph that still reads \"ToyBox *-- Toy\".\n# Paragraph.Range.Text includes the trailing paragraph-mark character\n# (CR), so trim it off before comparing.\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text.TrimEnd([char]13)\n    if ($t -eq \"ToyBox *-- Toy\") {\n        $target = $p.Range\n        break\n    }\n}\n\n# Step 3: rebuild that paragraph's runs: \"ToyBox\" stays wrapped in\n# spellStart/spellEnd proofing marks, followed by new runs for the\n# relationship text, and finish with the relocated bookmark.\n$ooxml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>' +\n  '<w:p w:rsidR=\"00C055B4\" w:rsidRDefault=\"00C055B4\" w:rsidP=\"00C055B4\">' +\n  '<w:proofErr w:type=\"spellStart\"/>' +\n  '<w:r><w:t>Toy</w:t></w:r>' +\n  '<w:r w:rsidR=\"00172387\"><w:t>B</w:t></w:r>' +\n  '<w:r><w:t>ox</w:t></w:r>' +\n  '<w:proofErr w:type=\"spellEnd\"/>' +\n  '<w:r><w:t xml:space=\"preserve\"> </w:t></w:r>' +\n  '<w:r><w:t>&quot;1&quot; *-- &quot;many&quot;</w:t></w:r>' +\n  '<w:r><w:t xml:space=\"preserve\"> Toy</w:t></w:r>' +\n  '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/>' +\n  '<w:bookmarkEnd w:id=\"0\"/>' +\n  '</w:p>' +\n  '</w:body></w:document>' +\n  '</pkg:xmlData></pkg:part></pkg:package>'\n\n$target.InsertXML($ooxml)\n"}
